{"js": "// Add yellow highlighting to the \"Open ledger\" bullet item (both the\n// paragraph mark and the run text), per the commit:\n// \"implementing a simple reusable form. Adding basic behavior to the\n// new ledger action, not yet operational.\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find((p) => p.text.trim() === \"Open ledger\");\nif (!target) {\n  throw new Error('Paragraph \"Open ledger\" not found.');\n}\n\n// Setting highlight on the paragraph's Font applies it to both the\n// paragraph mark (w:pPr/w:rPr) and the run(s) that make up its text\n// (w:r/w:rPr), matching the target edit.\ntarget.font.highlightColor = \"Yellow\";\n\nawait context.sync();\n", "ps1": "# Add yellow highlighting to the \"Open ledger\" bullet item (both the\n# paragraph mark and the run text), per the commit:\n# \"implementing a simple reusable form. Adding basic behavior to the\n# new ledger action, not yet operational.\"\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"Open ledger\") {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw 'Paragraph \"Open ledger\" not found.'\n}\n\n# Going through Range.Font (rather than Range.HighlightColorIndex\n# directly) highlights the paragraph mark (w:pPr/w:rPr) as well as the\n# run(s) making up the paragraph's text (w:r/w:rPr), matching the\n# target edit.\n$target.Range.Font.HighlightColorIndex = \"wdYellow\"\n"}
